# Aliens and Pyramid review - move the "Play Aliens and Pyramid for free -
# Review" / meta-description blurb from directly under the title down to
# just above the closing "image prompt" paragraph, splitting it in two:
#   - a new bold heading-style line "Play Aliens and Pyramid for free - Review"
#   - the trailing sentence replaces the old italic "Prompt: ..." text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Grab a formatted copy of paragraph 2 ("Meta description: Read our
#    review ...") *before* anything else moves around, so the run/format
#    structure (leading empty run + bold run + plain run) can be reused
#    for the newly inserted paragraph near the end of the document.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaFormatted = $metaPara.Range.FormattedText

# ---------------------------------------------------------------------
# 2) Build the new paragraph just above the final ("Prompt: ...")
#    paragraph, by inserting a fresh paragraph after the paragraph that
#    reads "The game's gameplay structure is standard." (the last bullet
#    under "What we don't like"), resetting it to the Normal style so it
#    does not inherit the bullet-list formatting, then pasting in the
#    formatted text captured above.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$standardPara = $d.Paragraphs($count - 1)
$standardPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($count)
$newPara.Style = "Normal"
$newPara.Range.FormattedText = $metaFormatted

# Trim the pasted text down to just the bold heading line: drop the
# ": Read our review ..." tail (the plain-formatted run) ...
$tail = $newPara.Range.Duplicate
$tail.Find.Execute(": Read our review of Aliens and Pyramid and play this unique slot game for free.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Delete()

# ... and rename "Meta description" (still bold) to the desired heading text.
$heading = $newPara.Range.Duplicate
$heading.Find.Execute("Meta description", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Play Aliens and Pyramid for free - Review", 2)

# ---------------------------------------------------------------------
# 3) Remove the original "Meta description" paragraph entirely (it now
#    lives, in split form, at the bottom of the document).
# ---------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------
# 4) Swap the final paragraph's italic "Prompt: ..." image-generation
#    text for the review blurb that used to follow "Meta description".
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Find.Execute( `
    "Prompt: Create a feature image for Aliens and Pyramids, the online slot game. The image should be in cartoon style and feature a happy Maya warrior with glasses. The image should be eye-catching and entice players to try out the game. The Maya warrior in the image should hold a golden key to depict the theme of the game, which is based on aliens and ancient Egypt. Use bright colors and be creative in designing the image.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our review of Aliens and Pyramid and play this unique slot game for free.", 2)
